$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New last row: "item" / "laptop" (written first, matching original authoring order)
$ws.Range("A13").Value = "item"
$ws.Range("B13").Value = "laptop"

# New personal/address test-data rows (5-12), replacing the old
# chusername/chpassword/pusername/ppassword block
$ws.Range("A5").Value = "firstName"
$ws.Range("B5").Value = "gael"

$ws.Range("A6").Value = "lastName"
$ws.Range("B6").Value = "gonzalez"

$ws.Range("A7").Value = "phoneNumber"
$ws.Range("B7").Value = "n1235454"
$ws.Range("B7").HorizontalAlignment = -4131

$ws.Range("A8").Value = "country"
$ws.Range("B8").Value = "mexico"

$ws.Range("A9").Value = "city"
$ws.Range("B9").Value = "ensenada"

$ws.Range("A10").Value = "address"
$ws.Range("B10").Value = "cercas"

$ws.Range("A11").Value = "state"
$ws.Range("B11").Value = "baja"

$ws.Range("A12").Value = "postalCode"
$ws.Range("B12").Value = "c22880"
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("B12").VerticalAlignment = -4160

# Row 1: username value changed last
$ws.Range("B1").Value = "julzzz1"

$ws.Range("C4").Select()
